# cyberdeck_bom.xlsx - "BOM" sheet
#
# 1) Analog stick (row 21) swapped for a 3DS slide pad, remarks updated.
# 2) A new row is inserted right after it (new row 22) for the FPC
#    connector used to attach that slide pad; every row from the old
#    row 22 onward shifts down by one (old 22->23, ... old 37->38).
# 3) dimension / autoFilter / the hidden _FilterDatabase name all grow
#    from ...G37 to ...G38 to track the now-38-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# --- 1) Update the analog-stick row --------------------------------------
$ws.Range("C21").Value = "3DSスライドパッド"
$ws.Range("G21").Value = "カーソル操作用 / FPC 6ピン 0.5mmピッチ 0.3mm厚 / 別途購入（AliExpress or 中古）"

# --- 2) Insert a new row at 22, pushing everything below down by one -----
$ws.Rows.Item(22).Insert()

# --- 3) Populate the new FPC-connector row --------------------------------
$ws.Range("A22").Value = "入力"
$ws.Range("B22").Value = "FPCコネクタ（スティック用）"
$ws.Range("C22").Value = "AFC01-S06FCA-00"
$ws.Range("D22").Value = "AFC01-S06FCA-00"
$ws.Range("E22").Value = "C262655"
# Leading apostrophe forces this numeric-looking quantity to be stored as
# text, matching every other "数量" cell in the sheet (all text, even "1").
$ws.Range("F22").Value = "'1"
$ws.Range("G22").Value = "6P 0.5mmピッチ Bottom Contact ZIF / 3DSスライドパッド接続用 / JLCPCB在庫52426個確認済み"

# Copy the formatting (fill/border/alignment) of row 21 onto the new row 22
# so it matches the surrounding "入力" block rather than the style that
# used to belong to row 22 ("オーディオ" block). Doing this *after* the
# value writes also re-normalises F22's style (the text-coercion above
# nudges a cell's style when it's written), so every cell in row 22 ends
# up on exactly the same style as row 21.
$ws.Range("A21:G21").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows.Item(22).RowHeight = 18

# --- 4) Grow the autofilter range from G37 to G38 -------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:G38").AutoFilter()

# ...and repoint the hidden _FilterDatabase defined name at the same range
# (inserting a row does not do this automatically).
$wb.Names.Item("BOM!_FilterDatabase").RefersTo = "='BOM'!`$A`$1:`$G`$38"
